# The commit swaps the content of the two DrawingML theme parts that ship
# with this deck: ppt/theme/theme1.xml (an unused "Office Theme" palette,
# only wired up via the Notes Master's relationship) and ppt/theme/theme2.xml
# (the "Integral"/"Red Violet" palette that is actually applied to the
# presentation via the Slide Master). After the commit, the presentation's
# live design uses the plain "Office" colour palette instead of the
# Red Violet one.
#
# The PowerPoint object model exposes exactly one editable theme colour
# scheme for the whole deck - SlideMaster.Theme.ThemeColorScheme (the same
# 12-slot scheme also reachable from NotesMaster/HandoutMaster/Slides, all
# of which resolve to the single live theme that is persisted as
# ppt/theme/theme2.xml). Re-pointing every slot to the stock "Office" RGB
# values reproduces the user-visible effect of the commit: the deck's
# design switches from the pink/violet "Integral" look to the standard
# blue/orange "Office" look.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

function RgbOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in the standard 12-slot COM order:
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $tcs.Item($i + 1).RGB = RgbOle $officeColors[$i]
}
